$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  2 = @{ C=1553; D=968 }
  3 = @{ C=1755; D=1542.5 }
  4 = @{ C=3851; D=3746 }
  5 = @{ C=4258; D=4210.5 }
  6 = @{ C=3197; D=3148.5 }
  7 = @{ C=1342; D=1459.5 }
  8 = @{ C=1876; D=1909 }
  9 = @{ C=426; D=963.5 }
  10 = @{ C=1331; D=1451.5 }
  11 = @{ C=3105; D=3026.5 }
  12 = @{ C=2798; D=2882 }
  13 = @{ C=639; D=429 }
  14 = @{ C=942; D=1010 }
  15 = @{ C=3377; D=3416 }
  16 = @{ C=1045; D=1070.5 }
  17 = @{ C=1329; D=1449.5 }
  18 = @{ C=1672; D=1488.5 }
  19 = @{ C=1928; D=1848 }
  20 = @{ C=3838; D=3708 }
  21 = @{ C=1022; D=1019 }
  22 = @{ C=1631; D=1043 }
  23 = @{ C=2432; D=2453.5 }
  24 = @{ C=966; D=927 }
  25 = @{ C=3241; D=3144 }
  26 = @{ C=4542; D=4491 }
  27 = @{ C=3608; D=3519 }
  28 = @{ C=1355; D=1400 }
  29 = @{ C=2037; D=2016 }
  30 = @{ C=3695; D=3699.5 }
  31 = @{ C=1619; D=1026 }
  32 = @{ C=1643; D=1482 }
  33 = @{ C=410; D=944 }
  34 = @{ C=889; D=877.5 }
  35 = @{ C=1750; D=1510.5 }
  36 = @{ C=3640; D=3786 }
  37 = @{ C=1706; D=1537.5 }
  38 = @{ C=1427; D=1470 }
  39 = @{ C=104; D=447.5 }
  40 = @{ C=1181; D=1063.5 }
  41 = @{ C=3435; D=3418 }
  42 = @{ C=355; D=302 }
  43 = @{ C=2521; D=2373.5 }
  44 = @{ C=1023; D=1020 }
  45 = @{ C=754; D=634.5 }
  46 = @{ C=794; D=697 }
  47 = @{ C=1200; D=1234.5 }
  48 = @{ C=611; D=381.5 }
  49 = @{ C=5135; D=5130 }
  50 = @{ C=4518; D=4522.5 }
  51 = @{ C=1373; D=1461.5 }
  52 = @{ C=2322; D=2268 }
  53 = @{ C=2356; D=1508.5 }
  54 = @{ C=2962; D=2846.5 }
  55 = @{ C=2047; D=1946 }
  56 = @{ C=601; D=357 }
  57 = @{ C=755; D=635.5 }
  58 = @{ C=2906; D=2876 }
  59 = @{ C=1694; D=1082.5 }
  60 = @{ C=2480; D=2807.5 }
  61 = @{ C=1654; D=1685 }
  62 = @{ C=2075; D=2045 }
  63 = @{ C=660; D=408 }
  64 = @{ C=566; D=566 }
  65 = @{ C=1011; D=1055.5 }
  66 = @{ C=4976; D=4975 }
  67 = @{ C=1483; D=1083 }
  68 = @{ C=802; D=670 }
  69 = @{ C=861; D=1783 }
  70 = @{ C=2585; D=2517 }
  71 = @{ C=3293; D=3238 }
  72 = @{ C=1315; D=1486 }
  73 = @{ C=1784; D=1602.5 }
  74 = @{ C=555; D=292 }
  75 = @{ C=102; D=448.5 }
  76 = @{ C=663; D=410.5 }
  77 = @{ C=2087; D=2064.5 }
  78 = @{ C=1613; D=1039 }
  79 = @{ C=1762; D=1662.5 }
  80 = @{ C=593; D=350 }
  81 = @{ C=777; D=672 }
  82 = @{ C=2043; D=2020 }
  83 = @{ C=2364; D=2308 }
  84 = @{ C=3840; D=3710 }
  85 = @{ C=87; D=468 }
  86 = @{ C=2225; D=2185 }
  87 = @{ C=4190; D=4208.5 }
  88 = @{ C=2852; D=2921 }
  89 = @{ C=784; D=664.5 }
  90 = @{ C=1104; D=1114.5 }
  91 = @{ C=705; D=448 }
  92 = @{ C=2691; D=2674.5 }
  93 = @{ C=2439; D=2407.5 }
  94 = @{ C=1701; D=1470 }
  95 = @{ C=3125; D=3066 }
  96 = @{ C=3101; D=3009.5 }
  97 = @{ C=3474; D=3679.5 }
  98 = @{ C=1768; D=1668.5 }
  99 = @{ C=1008; D=1014.5 }
  100 = @{ C=1440; D=1539.5 }
  101 = @{ C=1473; D=1558 }
  102 = @{ C=1711; D=1623 }
  103 = @{ C=2900; D=2846.5 }
  104 = @{ C=1800; D=1833.5 }
  105 = @{ C=3919; D=3741 }
  106 = @{ C=2562; D=2504 }
  107 = @{ C=878; D=600 }
  108 = @{ C=4567; D=4537 }
  109 = @{ C=1933; D=1855.5 }
  110 = @{ C=2156; D=2133 }
  111 = @{ C=2897; D=2809 }
  112 = @{ C=3134; D=3071.5 }
  113 = @{ C=1991.801801801802 }
}

foreach ($row in $updates.Keys) {
  $u = $updates[$row]
  if ($u.ContainsKey("C")) { $ws.Cells.Item($row, 3).Value = $u.C }
  if ($u.ContainsKey("D")) { $ws.Cells.Item($row, 4).Value = $u.D }
}
